# fixed initial bugs from beta test
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The "composite_id" row (originally row 2) moves to the bottom of the table
# (after "file_type", which was row 24, so it becomes row 24 once composite_id
# is removed from row 2).
$compositeField = $ws.Range("A2").Text
$compositeDef   = $ws.Range("B2").Text

# Remove row 2 entirely, shifting rows 3:24 up to 2:23.
$ws.Range("A2:B2").Delete(-4162) | Out-Null   # xlShiftUp

# Re-add composite_id as the new last row (row 24).
$ws.Range("A24").Value = $compositeField
$ws.Range("B24").Value = $compositeDef

# Several field definitions were clarified / updated during the beta-test fixes.
$ws.Range("B2").Value  = "university or organization of the uploading member (string)"
$ws.Range("B7").Value  = "current medical diagnosis or condition of the participant (string)"
$ws.Range("B8").Value  = "participant's smoking habits (e.g., yes/no)"
$ws.Range("B9").Value  = "current medical treatments or therapies the participant is undergoing (string)"
$ws.Range("B14").Value = "unique id of the study that the participant is taking part in (string)"
$ws.Range("B15").Value = "the duration (e.g., in days or months) between this scan and a previous scan for the same participant (string)"
$ws.Range("B23").Value = "file extension of image uploaded (string)"

$ws.Range("B22").Select() | Out-Null
